$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "96.487.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.74%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.614.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +8.69%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "640.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.49"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.96%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.403"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.55%  "
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.01"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.609.49"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +8.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.96%  "
$ws.Range("E13").Value = "  +4.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.37"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.290.60"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +8.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "96.485.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000254"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.618.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +9.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.04"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +21.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.497"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +11.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "516.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000199"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.75"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +10.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "98.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +10.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +23.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.144"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.63%  "
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.183"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.994"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "30.53"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.569"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.57%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "577.36"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.88"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.56%  "
$ws.Range("E39").Value = "  +8.40%  "
$ws.Range("E40").Value = "  +3.77%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.923"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.73%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0433"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.75"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.71"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "23.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.52%  "
$ws.Range("E47").Value = "  +5.46%  "
$ws.Range("E48").Value = "  -2.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.20"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.85%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.36%  "
